# fattyAcidChains.xlsx update
#  - add a new "PSM" column (H) to the FAS sheet:
#      H2  header  -> "PSM"
#      H3:H163 (data rows only, skipping the blank separator rows
#               32/54/76/98/120/142) -> the oxlipid/galactolipid rule string
#  - move the active selection to H2 and set the normal-view zoom to 100%,
#    matching the saved UI state of the edited workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rule string shared by every data row in the new column. The worksheet is
# laid out as stacked blocks of fatty-acid-chain rows separated by a blank
# spacer row (32, 54, 76, 98, 120, 142) — fill each contiguous block so the
# spacer rows themselves are left untouched. Written first so it lands in
# the shared-string table ahead of the "PSM" header (matches the source
# sharedStrings.xml ordering).
$rule = ";O1;O2;O3;O4;X1;X2,X3,X4;X1,O1;O2,X1;X2,O1"

$ws.Range("H3:H31").Value = $rule
$ws.Range("H33:H53").Value = $rule
$ws.Range("H55:H75").Value = $rule
$ws.Range("H77:H97").Value = $rule
$ws.Range("H99:H119").Value = $rule
$ws.Range("H121:H141").Value = $rule
$ws.Range("H143:H163").Value = $rule

# New header for the added column.
$ws.Range("H2").Value = "PSM"

# Match the saved view state: selection on the new header cell, 100% zoom.
$ws.Range("H2").Select()
$excel.ActiveWindow.Zoom = 100
